$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 8.6
$ws.Range("H3").Value = 1.45
$ws.Range("K3").Value = 5.4
$ws.Range("M3").Value = 1.02
$ws.Range("O3").Value = 1.14
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 1.14
$ws.Range("R3").Value = 1.37
$ws.Range("T3").Value = 1.78
$ws.Range("U3").Value = 1.01
$ws.Range("X3").Value = 1000
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 1000
$ws.Range("AC3").Value = 1000
$ws.Range("AD3").Value = 1000
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 1000
$ws.Range("AG3").Value = 1000
$ws.Range("AH3").Value = 1000
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 1000
$ws.Range("AK3").Value = 1000
$ws.Range("AL3").Value = 1000
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 1000
$ws.Range("AO3").Value = 1000

# Row 5
$ws.Range("N5").Value = 3.7
$ws.Range("U5").Value = 2.14

# Row 18
$ws.Range("G18").Value = 1.54
$ws.Range("J18").Value = 4.1
$ws.Range("P18").Value = 1.84
$ws.Range("Q18").Value = 2.1
$ws.Range("U18").Value = 1.53
$ws.Range("AI18").Value = 210
$ws.Range("AL18").Value = 70
$ws.Range("AM18").Value = 330

# Row 30
$ws.Range("K30").Value = 3.35

# Row 31
$ws.Range("I31").Value = 2.5

# Row 33
$ws.Range("AF33").Value = 110
$ws.Range("AJ33").Value = 430
$ws.Range("AK33").Value = 200
